$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix date format from "04-12-2023" to "04/12/2023"
#    A5 and C5 (and B8:B10 via the same shared string) all use this text.
#    Force text storage (not an auto-converted date serial number) by
#    applying a text number format before assigning the value, then
#    clearing the format again so no residual style is left behind.
$dateCells = @("A5", "C5", "B8", "B9", "B10")
foreach ($addr in $dateCells) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = "04/12/2023"
    $rng.ClearFormats()
}

# 2. Fix header capitalization "Total Hours per day" -> "Total Hours per Day"
$ws.Range("J7").Value = "Total Hours per Day"

# 3. Remove the empty "Notes" column values in G8:G10 (clears the now-unused empty shared string)
$ws.Range("G8").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("G10").ClearContents()

# 4. Merge header cells
$ws.Range("A1:C1").Merge()
$ws.Range("A3:B3").Merge()
